$d = $word.ActiveDocument

function Set-Highlight($searchText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $searchText
    $find.Replacement.Text = $searchText
    $find.Replacement.Highlight = 1
    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, `
                  $find.Replacement.Text, 1) | Out-Null
}

# 1) "Create a microservice environment consists of minimum 4 microservices." --
#    split off the trailing period into its own (unhighlighted) run, highlight the rest.
Set-Highlight "Create a microservice environment consists of minimum 4 microservices"

# 2) "There would be a swagger link for all services." -- highlight in three pieces so
#    the run boundaries land the same way the original authoring session produced them.
Set-Highlight "There would "
Set-Highlight "be a "
Set-Highlight "swagger link for all services."

# 3) "Each service has their own database." -- split off the leading space (left
#    unhighlighted) and highlight the sentence itself.
Set-Highlight "Each service has their own database."

# 4) "There should be at least 2 roles in the application." -- highlight whole sentence.
Set-Highlight "There should be at least 2 roles in the application."
